$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook is an auto-generated test-case report (UC005 - Listar Empenhos
# Pendentes). The v1.2.3 -> v1.2.4 fix swaps content that had landed in the
# wrong slots:
#   - TC1's 2nd step result ("success" message) and TC6's 2nd step result
#     ("MSG207 - unexpected error" message) were swapped.
#   - TC3's 1st step (action + result: "realizar o empenho" / "Apresenta a
#     tela de Registrar Empenho") and TC4's 1st step (action + result:
#     "atribuir/desatribuir" / "Atualiza a lista de registros") were swapped.

# --- Swap TC1 step-2 result (D11) with TC6 step-2 result (D48) ---
$tc1_step2_result = $ws.Range("D11").Value()
$tc6_step2_result = $ws.Range("D48").Value()
$ws.Range("D11").Value = $tc6_step2_result
$ws.Range("D48").Value = $tc1_step2_result

# --- Swap TC3 step-1 (B26/D26) with TC4 step-1 (B33/D33) ---
$tc3_step1_action = $ws.Range("B26").Value()
$tc3_step1_result = $ws.Range("D26").Value()
$tc4_step1_action = $ws.Range("B33").Value()
$tc4_step1_result = $ws.Range("D33").Value()

$ws.Range("B26").Value = $tc4_step1_action
$ws.Range("D26").Value = $tc4_step1_result
$ws.Range("B33").Value = $tc3_step1_action
$ws.Range("D33").Value = $tc3_step1_result
